$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-47 (Step, X, Y) following the DDA/Bresenham line algorithm
# Step = 0..45, X = Step, Y = round(Step * 29/45)
$data = @(
    @(0,0,0),
    @(1,1,1),
    @(2,2,1),
    @(3,3,2),
    @(4,4,3),
    @(5,5,3),
    @(6,6,4),
    @(7,7,5),
    @(8,8,5),
    @(9,9,6),
    @(10,10,6),
    @(11,11,7),
    @(12,12,8),
    @(13,13,8),
    @(14,14,9),
    @(15,15,10),
    @(16,16,10),
    @(17,17,11),
    @(18,18,12),
    @(19,19,12),
    @(20,20,13),
    @(21,21,14),
    @(22,22,14),
    @(23,23,15),
    @(24,24,15),
    @(25,25,16),
    @(26,26,17),
    @(27,27,17),
    @(28,28,18),
    @(29,29,19),
    @(30,30,19),
    @(31,31,20),
    @(32,32,21),
    @(33,33,21),
    @(34,34,22),
    @(35,35,23),
    @(36,36,23),
    @(37,37,24),
    @(38,38,24),
    @(39,39,25),
    @(40,40,26),
    @(41,41,26),
    @(42,42,27),
    @(43,43,28),
    @(44,44,28),
    @(45,45,29)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $triple = $data[$i]
    $ws.Cells.Item($row, 1).Value = $triple[0]
    $ws.Cells.Item($row, 2).Value = $triple[1]
    $ws.Cells.Item($row, 3).Value = $triple[2]
}
